$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +41.12%  "
$ws.Range("E6").Value = "  -4.27%  "
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +12.86%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  +9.08%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("E14").Value = "  -7.69%  "
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  +13.66%  "
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("E20").Value = "  +8.40%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  +9.90%  "
$ws.Range("E23").Value = "  +43.42%  "
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("E25").Value = "  -5.76%  "
$ws.Range("E26").Value = "  +12.24%  "
$ws.Range("E27").Value = "  -10.02%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("E42").Value = "  +3.82%  "
$ws.Range("E43").Value = "  +3.76%  "
$ws.Range("E44").Value = "  +9.96%  "
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("E49").Value = "  +14.55%  "
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("E51").Value = "  +1.94%  "
